$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# bch_usdt row (row 6) — add binance exchange limits for bch
$ws.Range("B6").Value = 0.01
$ws.Range("D6").Value = 0.00001
$ws.Range("F6").Value = 0.00001
$ws.Range("H6").Value = 10

# bch_btc row (row 10) — add binance exchange limits for bch
$ws.Range("B10").Value = 0.000001
$ws.Range("D10").Value = 0.001
$ws.Range("F10").Value = 0.001
$ws.Range("H10").Value = 0.001

# Update the active cell/selection to H10
$ws.Range("H10").Select()
